# Updated cryptos list (prices + 1h volume deltas) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''22.440.03'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '''1.567.04'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").Value = '''288.26'
$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("D7").Value = '''0.3724'
$ws.Range("E7").Value = '  +0.79%  '

$ws.Range("E8").Value = '  -3.34%  '

$ws.Range("D9").Value = '''0.3318'
$ws.Range("E9").Value = '  -2.04%  '

$ws.Range("D10").Value = '''0.07468'
$ws.Range("E10").Value = '  -1.22%  '

$ws.Range("D11").Value = '''1.130'
$ws.Range("E11").Value = '  -1.60%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").Value = '''20.74'
$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("D14").Value = '''5.952'
$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("D15").Value = '''6.899'
$ws.Range("E15").Value = '  -1.38%  '

$ws.Range("D16").Value = '''1.564.87'
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").Value = '''0.00001112'
$ws.Range("E17").Value = '  -0.88%  '

$ws.Range("D18").Value = '''0.06764'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = '''87.85'

$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").Value = '''6.349'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '''16.39'
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").Value = '''12.07'
$ws.Range("E23").Value = '  -0.99%  '

$ws.Range("D24").Value = '''22.428.49'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = '''2.387'
$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").Value = '''2.560'
$ws.Range("E26").Value = '  -4.15%  '

$ws.Range("D27").Value = '''153.70'
$ws.Range("E27").Value = '  +2.81%  '

$ws.Range("D28").Value = '''19.63'
$ws.Range("E28").Value = '  -1.87%  '

$ws.Range("D29").Value = '''5.016'
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("D30").Value = '''123.96'
$ws.Range("E30").Value = '  -0.97%  '

$ws.Range("D31").Value = '''1.741.30'
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("D32").Value = '''1.051'
$ws.Range("E32").Value = '  -1.43%  '

$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").Value = '''6.114'
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("D35").Value = '''9.624'
$ws.Range("E35").Value = '  -2.10%  '

$ws.Range("D36").Value = '''0.08298'
$ws.Range("E36").Value = '  -1.12%  '

$ws.Range("D37").Value = '''0.02453'
$ws.Range("E37").Value = '  -1.01%  '

$ws.Range("D38").Value = '''0.2270'
$ws.Range("E38").Value = '  -1.38%  '

$ws.Range("D39").Value = '''0.06376'
$ws.Range("E39").Value = '  -2.67%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.287'
$ws.Range("E40").Value = '  -4.61%  '

$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = '''5.347'
$ws.Range("E41").Value = '  -1.72%  '

$ws.Range("D42").Value = '''0.6265'
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").Value = '''11.22'
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").Value = '''1.002'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").Value = '''13.82'
$ws.Range("E45").Value = '  -1.69%  '

$ws.Range("D46").Value = '''0.6118'
$ws.Range("E46").Value = '  +4.10%  '

$ws.Range("D47").Value = '''3.779'
$ws.Range("E47").Value = '  -0.59%  '

$ws.Range("D48").Value = '''2.040'
$ws.Range("E48").Value = '  -1.53%  '

$ws.Range("D49").Value = '''125.42'
$ws.Range("E49").Value = '  -1.81%  '

$ws.Range("D50").Value = '''1.211'
$ws.Range("E50").Value = '  -2.50%  '

$ws.Range("D51").Value = '''0.07233'
$ws.Range("E51").Value = '  -0.99%  '
